# Redoing all the effects.
# Adds a new "Card Type" column (I) to the first worksheet with three
# category values (Anomaly / Room / Activator), matching the header style
# already used by the other column headers in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I1").Value = "Card Type"
$ws.Range("I2").Value = "Anomaly"
$ws.Range("I3").Value = "Room"
$ws.Range("I4").Value = "Activator"

# Size the new column to fit its (short) contents, same as the other
# bestFit columns on this sheet.
$ws.Columns.Item(9).ColumnWidth = 8.28

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("J7").Select() | Out-Null
